# Update with one column
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# Add the new "phonenumber" column to Sheet1
$ws1.Range("D1").Value = "phonenumber"
$ws1.Range("D2").Value = 123
$ws1.Range("D3").Value = 456
$ws1.Range("D4").Value = 789
[void]$ws1.Range("D5").Select()

# Rename Sheet3 to "check" and make it the active tab
$ws3.Name = "check"
$ws3.Activate()
